$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new week's data inserts two new rows right before the current row 280,
# pushing the existing rows 280-379 down to 282-381 (dimension grows from
# A1:R379 to A1:R381).
$ws.Rows("280:281").Insert()

# Duplicate the (now shifted) old row 280/281 formatting+values into the two
# freshly inserted blank rows, then overwrite with the new week's figures.
$ws.Range("A282:R282").Copy($ws.Range("A280:R280"))
$ws.Range("A283:R283").Copy($ws.Range("A281:R281"))

# New row 280: Zafiro rojo, 2022-09-21
$ws.Cells.Item(280, 4).Value = 44825
$ws.Cells.Item(280, 11).Value = 18000
$ws.Cells.Item(280, 12).Value = 19000
$ws.Cells.Item(280, 13).Value = 18500
$ws.Cells.Item(280, 16).Value = 1233

# New row 281: Zafiro verde, 2022-09-21
$ws.Cells.Item(281, 4).Value = 44825
$ws.Cells.Item(281, 11).Value = 15000
$ws.Cells.Item(281, 12).Value = 16000
$ws.Cells.Item(281, 13).Value = 15500
$ws.Cells.Item(281, 16).Value = 1033
